# Generate Report for Handback
# - Flip "Ready for handoff" status to "Handed back: in sync with en-US"
#   on the Overview sheet (both language columns) and on each language
#   sheet's Status column.
# - Refresh the "Latest Handback DateTime" timestamp for each language.
# - Clear the stale "Error Detail" message now that the handback is current.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-08-18 00:46:39"
$zhcn.Range("P2").Value = ""

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-08-18 00:46:47"
$dede.Range("P2").Value = ""
